$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Dane" (first sheet) - renumber the id column (A2:A7) down by one
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4
$ws1.Range("A6").Value = 5
$ws1.Range("A7").Value = 6

# ---------------------------------------------------------------------------
# Sheet "ocena" (second sheet) - add real score/rating data + number format
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 3.75
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 1

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 4.25
$ws2.Range("C3").Value = 6
$ws2.Range("D3").Value = 2

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 5
$ws2.Range("D4").Value = 3

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 2.25
$ws2.Range("C5").Value = 7
$ws2.Range("D5").Value = 4

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = 1.75
$ws2.Range("C6").Value = 10
$ws2.Range("D6").Value = 5

$ws2.Range("A7").Value = 6
$ws2.Range("B7").Value = 2
$ws2.Range("C7").Value = 5
$ws2.Range("D7").Value = 6

# Column B ("score") gets a 2-decimal number format
$ws2.Range("B1:B7").NumberFormat = "0.00"

# Give the "ocena" sheet a page setup (A4, portrait) like the "Dane" sheet
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Window / selection state
# ---------------------------------------------------------------------------

# "Dane": scroll so row 5 is at the top, select A8
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A8").Select()

# "ocena": zoom to 235%, select D5, and leave it as the active/tab-selected sheet
$ws2.Activate()
$excel.ActiveWindow.Zoom = 235
$ws2.Range("D5").Select()
